$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")
$ws.Range("B1").Value = "test"
